# Updates odds values on Sheet1 as per the commit "Atualizando o arquivo XLSX"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 2.63
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("X2").Value = 7
$ws.Range("AA2").Value = 19
$ws.Range("AW2").Value = 6.5

# Row 3
$ws.Range("BC3").Value = 151

# Row 5
$ws.Range("G5").Value = 4.2
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 1.73
$ws.Range("L5").Value = 2.3
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
$ws.Range("U5").Value = 1.67
$ws.Range("V5").Value = 2.1
$ws.Range("X5").Value = 23
$ws.Range("AA5").Value = 34
$ws.Range("AG5").Value = 8.5
$ws.Range("AM5").Value = 151
$ws.Range("AN5").Value = 6.5
$ws.Range("AO5").Value = 23
$ws.Range("AU5").Value = 7.5
$ws.Range("AY5").Value = 17
$ws.Range("BA5").Value = 41
$ws.Range("BB5").Value = 101

# Row 6
$ws.Range("G6").Value = 1.36
$ws.Range("H6").Value = 4.75
$ws.Range("I6").Value = 8.5
$ws.Range("W6").Value = 6.5
$ws.Range("Y6").Value = 8.5
$ws.Range("AG6").Value = 19
$ws.Range("AI6").Value = 23
$ws.Range("AJ6").Value = 101
$ws.Range("AW6").Value = 9
$ws.Range("BA6").Value = 201

# Row 8
$ws.Range("G8").Value = 2.63
$ws.Range("H8").Value = 2.8
$ws.Range("J8").Value = 3.5
$ws.Range("M8").Value = 1.13
$ws.Range("N8").Value = 6
$ws.Range("S8").Value = 1.62
$ws.Range("T8").Value = 2.2
$ws.Range("Z8").Value = 26
$ws.Range("AC8").Value = 6
$ws.Range("AD8").Value = 5.5
$ws.Range("AO8").Value = 17
$ws.Range("AP8").Value = 34
$ws.Range("AT8").Value = 2.2
$ws.Range("AU8").Value = 9.5
